$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing cell styles for the data range so that writing text-like
# numeric strings (e.g. "0.9993") does not get auto-converted to numbers by Excel,
# and does not leave a residual style/number-format change behind.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '24.439.10'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '1.685.52'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '315.35'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").Value = '0.3881'
$ws.Range("E7").Value = '  -1.31%  '
$ws.Range("D8").Value = '0.4019'
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").Value = '1.482'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").Value = '0.9991'
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").Value = '52.43'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("D12").Value = '0.08752'
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").Value = '25.94'
$ws.Range("E13").Value = '  +10.81%  '
$ws.Range("D14").Value = '7.483'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("D15").Value = '8.058'
$ws.Range("E15").Value = '  +0.46%  '
$ws.Range("D16").Value = '0.00001346'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").Value = '1.679.14'
$ws.Range("E17").Value = '  -2.18%  '
$ws.Range("D18").Value = '97.76'
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("D19").Value = '0.07280'
$ws.Range("E19").Value = '  +3.43%  '
$ws.Range("D20").Value = '19.82'
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").Value = '7.256'
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").Value = '14.17'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").Value = '24.426.11'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '3.018'
$ws.Range("E25").Value = '  -5.85%  '
$ws.Range("D26").Value = '2.341'
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("D27").Value = '22.56'
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = '167.34'
$ws.Range("E28").Value = '  +3.26%  '
$ws.Range("D29").Value = '8.617'
$ws.Range("E29").Value = '  +10.18%  '
$ws.Range("E30").Value = '  +3.63%  '
$ws.Range("D31").Value = '138.24'
$ws.Range("E31").Value = '  +1.24%  '
$ws.Range("D32").Value = '1.864.23'
$ws.Range("D33").Value = '0.08738'
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = '7.306'
$ws.Range("E34").Value = '  +1.75%  '
$ws.Range("D35").Value = '1.044'
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '2.043'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.03051'
$ws.Range("E37").Value = '  +11.05%  '
$ws.Range("D38").Value = '0.2773'
$ws.Range("E38").Value = '  +1.18%  '
$ws.Range("D39").Value = '10.78'
$ws.Range("E39").Value = '  -3.75%  '
$ws.Range("D40").Value = '0.09128'
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("D41").Value = '0.7994'
$ws.Range("E41").Value = '  +4.18%  '
$ws.Range("D42").Value = '14.09'
$ws.Range("E42").Value = '  -2.02%  '
$ws.Range("D43").Value = '1.470'
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("D44").Value = '17.68'
$ws.Range("E44").Value = '  +12.89%  '
$ws.Range("D45").Value = '2.623'
$ws.Range("E45").Value = '  +1.93%  '
$ws.Range("D46").Value = '0.7232'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").Value = '4.262'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").Value = '1.421'
$ws.Range("E48").Value = '  +8.11%  '
$ws.Range("D49").Value = '0.9986'
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("D50").Value = '138.97'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").Value = '0.08070'
$ws.Range("E51").Value = '  +0.95%  '

# Restore original styling/number format so only cell contents changed.
$dataRange.Style = $origStyle

